$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the five existing "sac7x@grr.la" cells to the new "sac8x@grr.la" values.
$ws.Range("C1").Value = "sac80@grr.la"
$ws.Range("C2").Value = "sac81@grr.la"
$ws.Range("C16").Value = "sac83@grr.la"
$ws.Range("C17").Value = "sac84@grr.la"
$ws.Range("C18").Value = "sac85@grr.la"

# Add the new G2 cell with the sixth new email address, styled like the other
# hyperlink cells.
$ws.Range("G2").Value = "sac82@grr.la"
$ws.Range("G2").Style = "Hyperlink"

# Rebuild the hyperlinks collection: the underlying host does not support
# in-place updates of a hyperlink's Address (it only appends new entries), so
# clear everything and re-add each hyperlink with its correct target in the
# same order as before, including the unaffected ones (A3, A21).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:sac1@grr.la", "", "", "sac1@grr.la")
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:sac81@grr.la")
$ws.Hyperlinks.Add($ws.Range("C1"), "mailto:sac80@grr.la")
$ws.Hyperlinks.Add($ws.Range("C16"), "mailto:sac83@grr.la")
$ws.Hyperlinks.Add($ws.Range("C17"), "mailto:sac84@grr.la")
$ws.Hyperlinks.Add($ws.Range("C18"), "mailto:sac85@grr.la")
$ws.Hyperlinks.Add($ws.Range("A21"), "mailto:sac50@grr.la")
$ws.Hyperlinks.Add($ws.Range("G2"), "mailto:sac82@grr.la")

# Update the view: selection moves to C19 and the frozen/top-left cell resets
# to the default (A1) instead of A3.
$ws.Range("C19").Select()
